# Updated cryptos list on Sat Sep  2 04:36:54 UTC 2023 with GitHub Actions
#
# Refreshes the live snapshot of the crypto table on Sheet1: column D
# (Price) and column E (Volume(1h)) are updated to the latest scraped
# values for every coin row. Rows 41/42 additionally swap ranking order
# (PaxDollar now ranks above mCoin), so their Coin name (B) and Link (C)
# columns are updated as well.
#
# Price values are written with a leading text qualifier where needed
# so that figures such as "1.002" or "15.70" stay as the literal text
# scraped from the site instead of being re-interpreted as numbers
# (matching how these cells were already stored as text in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.796.38'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').Value = '1.635.02'
$ws.Range('E3').Value = '  -1.47%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('E6').Value = '  -2.26%  '
$ws.Range('D7').Value = '''1.002'
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = '''0.06386'
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('D10').Value = '''19.68'
$ws.Range('E10').Value = '  -1.62%  '
$ws.Range('D11').Value = '''0.07696'
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('D12').Value = '''4.254'
$ws.Range('E12').Value = '  -0.98%  '
$ws.Range('D13').Value = '1.634.12'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').Value = '1.859.34'
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('D15').Value = '''0.5458'
$ws.Range('E15').Value = '  -1.64%  '
$ws.Range('D16').Value = '0.0₅7939'
$ws.Range('E16').Value = '  -1.56%  '
$ws.Range('D17').Value = '''63.58'
$ws.Range('E17').Value = '  -1.11%  '
$ws.Range('D18').Value = '25.837.23'
$ws.Range('E18').Value = '  -1.50%  '
$ws.Range('E19').Value = '  -0.28%  '
$ws.Range('E20').Value = '  -3.68%  '
$ws.Range('D21').Value = '''4.328'
$ws.Range('E21').Value = '  -2.31%  '
$ws.Range('D22').Value = '''9.948'
$ws.Range('E22').Value = '  -1.38%  '
$ws.Range('D23').Value = '''5.971'
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = '''1.916'
$ws.Range('E25').Value = '  +9.25%  '
$ws.Range('D26').Value = '''141.30'
$ws.Range('E26').Value = '  -2.72%  '
$ws.Range('D27').Value = '''0.1147'
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('D28').Value = '''15.70'
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('D29').Value = '''6.713'
$ws.Range('E29').Value = '  -4.04%  '
$ws.Range('D30').Value = '''0.05025'
$ws.Range('E30').Value = '  -3.62%  '
$ws.Range('D31').Value = '''1.242'
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('D32').Value = '''3.273'
$ws.Range('E32').Value = '  -2.53%  '
$ws.Range('D33').Value = '''3.189'
$ws.Range('E33').Value = '  -1.38%  '
$ws.Range('D34').Value = '''1.539'
$ws.Range('E34').Value = '  -2.30%  '
$ws.Range('D35').Value = '''2.353'
$ws.Range('E35').Value = '  -0.86%  '
$ws.Range('D36').Value = '1.175.18'
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').Value = '''0.8932'
$ws.Range('E37').Value = '  -4.10%  '
$ws.Range('D38').Value = '''2.608'
$ws.Range('E38').Value = '  -5.55%  '
$ws.Range('D39').Value = '''0.5602'
$ws.Range('E39').Value = '  -1.82%  '
$ws.Range('D40').Value = '''0.01561'
$ws.Range('E40').Value = '  -2.26%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '''1.002'
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('B42').Value = 'mCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D42').Value = '''2.545'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').Value = '''5.670'
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('D44').Value = '''0.8088'
$ws.Range('E44').Value = '  -3.61%  '
$ws.Range('D45').Value = '''99.55'
$ws.Range('E45').Value = '  -1.13%  '
$ws.Range('D46').Value = '1.771.64'
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('E47').Value = '  -0.41%  '
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('D49').Value = '''1.003'
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('E50').Value = '  -1.89%  '
$ws.Range('E51').Value = '  -0.48%  '
